$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row at position 17 (shifts old rows 17..58 down to 18..59),
#    copying the formatting of row 16 (the normal data-row style) into the
#    newly inserted row 17.
# ------------------------------------------------------------------
$ws.Rows("17:17").Insert()
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Fill the new row 17 with the new period (2507) for the same worker
#    as the rest of the table (CC 1048217270 - CRISTIAN DE JUESUS PALMA DE LA RANS).
# ------------------------------------------------------------------
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1048217270"
$ws.Range("D17").Value = "CRISTIAN DE JUESUS PALMA DE LA RANS"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1000000

# ------------------------------------------------------------------
# 3. Re-write the "Periodo Mora" column for the other period rows so that
#    they now run in descending order (2506 .. 2208) instead of ascending
#    (2208 .. 2506). These rows are now at 18..52 after the insert.
# ------------------------------------------------------------------
$periods = @(
    "2506","2505","2504","2503","2502","2501",
    "2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401",
    "2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301",
    "2212","2211","2210","2209","2208"
)

$row = 18
foreach ($p in $periods) {
    $ws.Cells.Item($row, 5).Value = $p
    $ws.Cells.Item($row, 6).Value = 40000
    $row++
}

# ------------------------------------------------------------------
# 4. The row that used to be the last data row (old row 17, with the bottom
#    border style) is now row 53; restore its original period/value
#    (2207 / 16000).
# ------------------------------------------------------------------
$ws.Range("E53").Value = "2207"
$ws.Range("F53").Value = 16000

# ------------------------------------------------------------------
# 5. Update the summary figures: total overdue value grew by 40000 (the new
#    period) and the period count grew by one.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 1486666
$ws.Range("F13").Value = 37
